# Cookie.xlsx proto-table update
#  - row 3 "용감한 쿠기" (typo) -> "용감한 쿠키", id 10001 -> 1001
#  - new row 4: id 1002, name "딸기맛 쿠키"
#  - selection moves to A5 (next empty row) after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row first so the brand-new shared string ("딸기맛 쿠키")
# is interned before the corrected-typo string, matching the order the
# strings end up in inside sharedStrings.xml.
$ws.Range("A4").Value = 1002
$ws.Range("B4").Value = "딸기맛 쿠키"

# Fix row 3: the numeric id and the misspelled cookie name.
$ws.Range("A3").Value = 1001
$ws.Range("B3").Value = "용감한 쿠키"

# Leave the selection where the author left it - the next blank row.
$ws.Range("A5").Select()
